$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(42, 8).Value = 7802.2
$ws.Cells.Item(42, 9).Value = 6007.3335
$ws.Cells.Item(42, 10).Value = 10494.5
$ws.Cells.Item(42, 11).Value = 18022.0005
$ws.Cells.Item(42, 12).Value = 31483.5
$ws.Cells.Item(42, 13).Value = -17792.0005
$ws.Cells.Item(42, 14).Value = -31943.5

$ws.Cells.Item(51, 8).Value = 13799.4
$ws.Cells.Item(51, 9).Value = 9998
$ws.Cells.Item(51, 10).Value = 14749.75
$ws.Cells.Item(51, 11).Value = 9998
$ws.Cells.Item(51, 12).Value = 14749.75
$ws.Cells.Item(51, 13).Value = -9514
$ws.Cells.Item(51, 14).Value = -15717.75

$ws.Cells.Item(86, 8).Value = 8598.444
$ws.Cells.Item(86, 9).Value = 8172.8076
$ws.Cells.Item(86, 10).Value = 9705.1
$ws.Cells.Item(86, 11).Value = 8172.8076
$ws.Cells.Item(86, 12).Value = 9705.1
$ws.Cells.Item(86, 13).Value = -7049.8076
$ws.Cells.Item(86, 14).Value = -11951.1

$ws.Cells.Item(89, 8).Value = 8598.444
$ws.Cells.Item(89, 9).Value = 8172.8076
$ws.Cells.Item(89, 10).Value = 9705.1
$ws.Cells.Item(89, 11).Value = 40864.038
$ws.Cells.Item(89, 12).Value = 48525.5
$ws.Cells.Item(89, 13).Value = -35248.038
$ws.Cells.Item(89, 14).Value = -59757.5

$ws.Cells.Item(116, 8).Value = 8603.759
$ws.Cells.Item(116, 9).Value = 10306
$ws.Cells.Item(116, 11).Value = 10306
$ws.Cells.Item(116, 13).Value = -6864

$ws.Cells.Item(119, 8).Value = 1500
$ws.Cells.Item(119, 10).Value = 1500
$ws.Cells.Item(119, 12).Value = 4500
$ws.Cells.Item(119, 14).Value = -14176

$ws.Cells.Item(129, 8).Value = 1553.8334
$ws.Cells.Item(129, 10).Value = 1399.8
$ws.Cells.Item(129, 12).Value = 4199.4
$ws.Cells.Item(129, 14).Value = -14199.4

$ws.Cells.Item(135, 8).Value = 93750510
$ws.Cells.Item(135, 9).Value = 35714840
$ws.Cells.Item(135, 11).Value = 321433560
$ws.Cells.Item(135, 13).Value = -321431025

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2101498.8
$ws.Cells.Item(2, 9).Value = 2451348.5
$ws.Cells.Item(2, 11).Value = 2451348.5
$ws.Cells.Item(2, 13).Value = -2451235.5

$ws.Cells.Item(7, 8).Value = 99000
$ws.Cells.Item(7, 10).Value = 99000
$ws.Cells.Item(7, 12).Value = 99000
$ws.Cells.Item(7, 14).Value = -99228

$ws.Cells.Item(16, 8).Value = 9998.4
$ws.Cells.Item(16, 10).Value = 10399.333
$ws.Cells.Item(16, 12).Value = 10399.333
$ws.Cells.Item(16, 14).Value = -10973.333

$ws.Cells.Item(45, 8).Value = 1406.2
$ws.Cells.Item(45, 9).Value = 1316.2307
$ws.Cells.Item(45, 11).Value = 1316.2307
$ws.Cells.Item(45, 13).Value = -939.2307000000001

$ws.Cells.Item(51, 8).Value = 20042
$ws.Cells.Item(51, 9).Value = 20042
$ws.Cells.Item(51, 11).Value = 20042
$ws.Cells.Item(51, 13).Value = -19286

$ws.Cells.Item(61, 8).Value = 166671330
$ws.Cells.Item(61, 9).Value = 166671330
$ws.Cells.Item(61, 11).Value = 166671330
$ws.Cells.Item(61, 13).Value = -166671118

$ws.Cells.Item(74, 8).Value = 58826930
$ws.Cells.Item(74, 9).Value = 71432200
$ws.Cells.Item(74, 10).Value = 2333
$ws.Cells.Item(74, 11).Value = 71432200
$ws.Cells.Item(74, 12).Value = 2333
$ws.Cells.Item(74, 13).Value = -71431326
$ws.Cells.Item(74, 14).Value = -4081

$ws.Cells.Item(77, 8).Value = 58826930
$ws.Cells.Item(77, 9).Value = 71432200
$ws.Cells.Item(77, 10).Value = 2333
$ws.Cells.Item(77, 11).Value = 357161000
$ws.Cells.Item(77, 12).Value = 11665
$ws.Cells.Item(77, 13).Value = -357156632
$ws.Cells.Item(77, 14).Value = -20401

$ws.Cells.Item(110, 8).Value = 35701.035
$ws.Cells.Item(110, 9).Value = 37894.355
$ws.Cells.Item(110, 11).Value = 37894.355
$ws.Cells.Item(110, 13).Value = -35849.355

$ws.Cells.Item(116, 8).Value = 2101498.8
$ws.Cells.Item(116, 9).Value = 2451348.5
$ws.Cells.Item(116, 11).Value = 2451348.5
$ws.Cells.Item(116, 13).Value = -2449054.5

$ws.Cells.Item(136, 8).Value = 166671330
$ws.Cells.Item(136, 9).Value = 166671330
$ws.Cells.Item(136, 11).Value = 500013990
$ws.Cells.Item(136, 13).Value = -500011440

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2101498.8
$ws.Cells.Item(3, 9).Value = 2451348.5
$ws.Cells.Item(3, 11).Value = 2451348.5
$ws.Cells.Item(3, 13).Value = -2451234.5

$ws.Cells.Item(26, 8).Value = 22221
$ws.Cells.Item(26, 9).Value = 22221
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 22221
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = -21929
$ws.Cells.Item(26, 14).Value = $null

$ws.Cells.Item(86, 8).Value = 4061.5833
$ws.Cells.Item(86, 9).Value = 4443.5557
$ws.Cells.Item(86, 10).Value = 2915.6667
$ws.Cells.Item(86, 11).Value = 4443.5557
$ws.Cells.Item(86, 12).Value = 2915.6667
$ws.Cells.Item(86, 13).Value = -3320.5557
$ws.Cells.Item(86, 14).Value = -5161.6667

$ws.Cells.Item(89, 8).Value = 4061.5833
$ws.Cells.Item(89, 9).Value = 4443.5557
$ws.Cells.Item(89, 10).Value = 2915.6667
$ws.Cells.Item(89, 11).Value = 22217.7785
$ws.Cells.Item(89, 12).Value = 14578.3335
$ws.Cells.Item(89, 13).Value = -16601.7785
$ws.Cells.Item(89, 14).Value = -25810.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(17, 8).Value = 500
$ws.Cells.Item(17, 9).Value = 500
$ws.Cells.Item(17, 11).Value = 500
$ws.Cells.Item(17, 13).Value = -326

$ws.Cells.Item(31, 8).Value = 13938.765
$ws.Cells.Item(31, 9).Value = 5137
$ws.Cells.Item(31, 10).Value = 20100
$ws.Cells.Item(31, 11).Value = 5137
$ws.Cells.Item(31, 12).Value = 20100
$ws.Cells.Item(31, 13).Value = -4842
$ws.Cells.Item(31, 14).Value = -20690

$ws.Cells.Item(34, 8).Value = 13938.765
$ws.Cells.Item(34, 9).Value = 5137
$ws.Cells.Item(34, 10).Value = 20100
$ws.Cells.Item(34, 11).Value = 5137
$ws.Cells.Item(34, 12).Value = 20100
$ws.Cells.Item(34, 13).Value = -4935
$ws.Cells.Item(34, 14).Value = -20504

$ws.Cells.Item(58, 8).Value = 26323232
$ws.Cells.Item(58, 9).Value = 29419848
$ws.Cells.Item(58, 10).Value = 2007
$ws.Cells.Item(58, 11).Value = 29419848
$ws.Cells.Item(58, 12).Value = 2007
$ws.Cells.Item(58, 13).Value = -29419645
$ws.Cells.Item(58, 14).Value = -2413

$ws.Cells.Item(94, 8).Value = 844
$ws.Cells.Item(94, 9).Value = 761.8333
$ws.Cells.Item(94, 10).Value = 967.25
$ws.Cells.Item(94, 11).Value = 761.8333
$ws.Cells.Item(94, 12).Value = 967.25
$ws.Cells.Item(94, 13).Value = -310.8333
$ws.Cells.Item(94, 14).Value = -1869.25

$ws.Cells.Item(136, 8).Value = 26323232
$ws.Cells.Item(136, 9).Value = 29419848
$ws.Cells.Item(136, 10).Value = 2007
$ws.Cells.Item(136, 11).Value = 88259544
$ws.Cells.Item(136, 12).Value = 6021
$ws.Cells.Item(136, 13).Value = -88256994
$ws.Cells.Item(136, 14).Value = -11121

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1356864.2
$ws.Cells.Item(4, 9).Value = 1252169.8
$ws.Cells.Item(4, 10).Value = 1431645.9
$ws.Cells.Item(4, 11).Value = 3756509.4
$ws.Cells.Item(4, 12).Value = 4294937.699999999
$ws.Cells.Item(4, 13).Value = -3756397.4
$ws.Cells.Item(4, 14).Value = -4295161.699999999

$ws.Cells.Item(23, 8).Value = 1129.625
$ws.Cells.Item(23, 9).Value = 1004.5
$ws.Cells.Item(23, 11).Value = 3013.5
$ws.Cells.Item(23, 13).Value = -2778.5

$ws.Cells.Item(46, 8).Value = 297
$ws.Cells.Item(46, 10).Value = 504
$ws.Cells.Item(46, 12).Value = 1512
$ws.Cells.Item(46, 14).Value = -1694

$ws.Cells.Item(56, 8).Value = 18818.965
$ws.Cells.Item(56, 9).Value = 18818.965
$ws.Cells.Item(56, 11).Value = 18818.965
$ws.Cells.Item(56, 13).Value = -18288.965

$ws.Cells.Item(68, 8).Value = 2945
$ws.Cells.Item(68, 9).Value = 720
$ws.Cells.Item(68, 11).Value = 2160
$ws.Cells.Item(68, 13).Value = -1349

$ws.Cells.Item(71, 8).Value = 2945
$ws.Cells.Item(71, 9).Value = 720
$ws.Cells.Item(71, 11).Value = 6480
$ws.Cells.Item(71, 13).Value = -2424

$ws.Cells.Item(131, 8).Value = 5683.3335
$ws.Cells.Item(131, 9).Value = 5625
$ws.Cells.Item(131, 11).Value = 16875
$ws.Cells.Item(131, 13).Value = -11835

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 680.1111
$ws.Cells.Item(2, 9).Value = 812
$ws.Cells.Item(2, 11).Value = 812
$ws.Cells.Item(2, 13).Value = -699

$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 14).Value = $null

$ws.Cells.Item(29, 8).Value = 16666
$ws.Cells.Item(29, 9).Value = 16666
$ws.Cells.Item(29, 11).Value = 16666
$ws.Cells.Item(29, 13).Value = -16376

$ws.Cells.Item(36, 8).Value = 12500
$ws.Cells.Item(36, 10).Value = 12500
$ws.Cells.Item(36, 12).Value = 12500
$ws.Cells.Item(36, 14).Value = -13470

$ws.Cells.Item(40, 8).Value = 23998
$ws.Cells.Item(40, 9).Value = 26000
$ws.Cells.Item(40, 10).Value = 22663.334
$ws.Cells.Item(40, 11).Value = 26000
$ws.Cells.Item(40, 12).Value = 22663.334
$ws.Cells.Item(40, 13).Value = -25849
$ws.Cells.Item(40, 14).Value = -22965.334

$ws.Cells.Item(97, 8).Value = 1034.6316
$ws.Cells.Item(97, 9).Value = 619.8182
$ws.Cells.Item(97, 11).Value = 619.8182
$ws.Cells.Item(97, 13).Value = -123.8182

$ws.Cells.Item(104, 8).Value = 38972
$ws.Cells.Item(104, 10).Value = 38972
$ws.Cells.Item(104, 12).Value = 38972
$ws.Cells.Item(104, 14).Value = -45960

$ws.Cells.Item(122, 8).Value = 89733.78999999999
$ws.Cells.Item(122, 9).Value = 134919.22
$ws.Cells.Item(122, 10).Value = 8400
$ws.Cells.Item(122, 11).Value = 404757.66
$ws.Cells.Item(122, 12).Value = 25200
$ws.Cells.Item(122, 13).Value = -402307.66
$ws.Cells.Item(122, 14).Value = -30100

$ws.Cells.Item(126, 8).Value = 4628.5
$ws.Cells.Item(126, 9).Value = 4954.7
$ws.Cells.Item(126, 10).Value = 2997.5
$ws.Cells.Item(126, 11).Value = 14864.1
$ws.Cells.Item(126, 12).Value = 8992.5
$ws.Cells.Item(126, 13).Value = -12394.1
$ws.Cells.Item(126, 14).Value = -13932.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(106, 8).Value = 11249
$ws.Cells.Item(106, 10).Value = 11249
$ws.Cells.Item(106, 12).Value = 11249
$ws.Cells.Item(106, 14).Value = -13773

$ws.Cells.Item(122, 8).Value = 6055.5884
$ws.Cells.Item(122, 9).Value = 6463
$ws.Cells.Item(122, 11).Value = 19389
$ws.Cells.Item(122, 13).Value = -16939

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 9999
$ws.Cells.Item(3, 10).Value = 9999
$ws.Cells.Item(3, 12).Value = 9999
$ws.Cells.Item(3, 14).Value = -10227

$ws.Cells.Item(25, 8).Value = 23500
$ws.Cells.Item(25, 10).Value = 23500
$ws.Cells.Item(25, 12).Value = 23500
$ws.Cells.Item(25, 14).Value = -24086

$ws.Cells.Item(136, 8).Value = 10641092
$ws.Cells.Item(136, 9).Value = 11907683
$ws.Cells.Item(136, 11).Value = 35723049
$ws.Cells.Item(136, 13).Value = -35720499
